# Apply updated crypto price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.755.72"
$ws.Range("E2").Value = "  +3.73%  "
$ws.Range("D3").Value = "'3.147.18"
$ws.Range("E3").Value = "  +2.54%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'578.35"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "'180.00"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D8").Value = "'3.146.07"
$ws.Range("E8").Value = "  +2.60%  "
$ws.Range("D9").Value = "'0.524"
$ws.Range("E9").Value = "  +2.83%  "
$ws.Range("D10").Value = "'6.52"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "'0.471"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  +1.54%  "
$ws.Range("E14").Value = "  +3.64%  "
$ws.Range("D15").Value = "'68.595.26"
$ws.Range("E15").Value = "  +3.53%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.122"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "'3.669.11"
$ws.Range("E17").Value = "  +2.49%  "
$ws.Range("D18").Value = "'7.15"
$ws.Range("E18").Value = "  +2.76%  "
$ws.Range("D19").Value = "'3.144.55"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("D20").Value = "'16.45"
$ws.Range("E20").Value = "  -2.91%  "
$ws.Range("D21").Value = "'489.98"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'0.698"
$ws.Range("E22").Value = "  +1.58%  "
$ws.Range("D23").Value = "'7.79"
$ws.Range("E23").Value = "  +0.57%  "
$ws.Range("D24").Value = "'84.09"
$ws.Range("E24").Value = "  +1.78%  "
$ws.Range("E25").Value = "  +6.19%  "
$ws.Range("D26").Value = "'13.00"
$ws.Range("E26").Value = "  +2.39%  "
$ws.Range("D27").Value = "'10.58"
$ws.Range("E27").Value = "  +3.99%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'8.14"
$ws.Range("E29").Value = "  +4.32%  "
$ws.Range("E30").Value = "  +4.25%  "
$ws.Range("E31").Value = "  +1.44%  "
$ws.Range("D32").Value = "'28.18"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("E33").Value = "  +0.96%  "
$ws.Range("E34").Value = "  +4.36%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'5.77"
$ws.Range("E36").Value = "  +3.11%  "
$ws.Range("D37").Value = "'47.99"
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("D38").Value = "'0.960"
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").Value = "'0.323"
$ws.Range("E39").Value = "  +7.38%  "
$ws.Range("E40").Value = "  +3.81%  "
$ws.Range("E41").Value = "  +3.23%  "
$ws.Range("D42").Value = "'49.25"
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'8.37"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "'400.60"
$ws.Range("E44").Value = "  +9.45%  "
$ws.Range("D45").Value = "'2.72"
$ws.Range("E45").Value = "  +7.20%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "'2.814.12"
$ws.Range("E46").Value = "  +1.20%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'27.55"
$ws.Range("E47").Value = "  +12.86%  "
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("D49").Value = "'135.17"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D51").Value = "'2.37"
$ws.Range("E51").Value = "  +9.81%  "
